$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text so values like "1.000" or "0.00001098" are not
# reinterpreted as numbers/dates by the COM value setter.
$ws.Range("D2:D51").NumberFormat = "@"

# Updated crypto price/volume snapshot (and three re-ranked rows: 26-28)
$ws.Range('D2').Value = '27.331.78'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').Value = '1.784.86'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '340.62'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '0.3963'
$ws.Range('E7').Value = '  +3.63%  '
$ws.Range('D8').Value = '0.3463'
$ws.Range('E8').Value = '  -1.93%  '
$ws.Range('D9').Value = '47.97'
$ws.Range('E9').Value = '  -3.83%  '
$ws.Range('D10').Value = '1.198'
$ws.Range('E10').Value = '  -3.07%  '
$ws.Range('D11').Value = '0.07518'
$ws.Range('E11').Value = '  -2.73%  '
$ws.Range('D12').Value = '0.9992'
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('D13').Value = '21.81'
$ws.Range('E13').Value = '  -2.63%  '
$ws.Range('D14').Value = '6.478'
$ws.Range('E14').Value = '  -2.10%  '
$ws.Range('D15').Value = '1.783.76'
$ws.Range('D16').Value = '7.116'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '0.00001098'
$ws.Range('E17').Value = '  -2.56%  '
$ws.Range('D18').Value = '0.06698'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').Value = '84.93'
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('D20').Value = '0.9998'
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').Value = '17.72'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').Value = '6.517'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').Value = '27.328.34'
$ws.Range('E23').Value = '  -0.74%  '
$ws.Range('D24').Value = '12.43'
$ws.Range('E24').Value = '  -5.44%  '
$ws.Range('D25').Value = '2.383'
$ws.Range('E25').Value = '  -4.08%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '21.27'
$ws.Range('E26').Value = '  -3.60%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.501'
$ws.Range('E27').Value = '  -6.53%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').Value = '1.460'
$ws.Range('E28').Value = '  -1.43%  '
$ws.Range('D29').Value = '157.32'
$ws.Range('E29').Value = '  +2.76%  '
$ws.Range('D30').Value = '1.985.91'
$ws.Range('E30').Value = '  -2.19%  '
$ws.Range('D31').Value = '136.40'
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('D32').Value = '4.033'
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('D33').Value = '5.988'
$ws.Range('E33').Value = '  -5.45%  '
$ws.Range('D34').Value = '0.08848'
$ws.Range('E34').Value = '  +0.58%  '
$ws.Range('D35').Value = '13.02'
$ws.Range('E35').Value = '  -6.73%  '
$ws.Range('D36').Value = '0.02453'
$ws.Range('E36').Value = '  +2.02%  '
$ws.Range('D37').Value = '1.618'
$ws.Range('E37').Value = '  -4.76%  '
$ws.Range('D38').Value = '5.421'
$ws.Range('E38').Value = '  -3.84%  '
$ws.Range('D39').Value = '0.06482'
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').Value = '0.6854'
$ws.Range('E40').Value = '  -2.52%  '
$ws.Range('D41').Value = '0.2217'
$ws.Range('E41').Value = '  -1.98%  '
$ws.Range('D42').Value = '1.255'
$ws.Range('E42').Value = '  -3.18%  '
$ws.Range('D43').Value = '8.382'
$ws.Range('E43').Value = '  -8.26%  '
$ws.Range('D44').Value = '14.54'
$ws.Range('E44').Value = '  -1.60%  '
$ws.Range('D45').Value = '0.9992'
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('D46').Value = '0.6406'
$ws.Range('E46').Value = '  -3.28%  '
$ws.Range('E47').Value = '  -1.46%  '
$ws.Range('D48').Value = '2.139'
$ws.Range('E48').Value = '  -2.33%  '
$ws.Range('D49').Value = '132.47'
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('D50').Value = '0.07146'
$ws.Range('E50').Value = '  -2.23%  '
$ws.Range('D51').Value = '79.43'
$ws.Range('E51').Value = '  -2.17%  '

# Restore the original (default) cell style on the Price column now that the
# text values are committed, so no stray number-format styling is introduced.
$ws.Range("D2:D51").Style = "Normal"
